# tried and added xgboost and random forest
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 13: random forest approach
$ws.Range("C13").Value = "random forest"
$ws.Range("D13").Value = "do"
$ws.Range("E13").Value = "NA"
$ws.Range("F13").Value = "do"
$ws.Range("H13").Value = "RandomForestRegressor(max_depth=5, random_state=12, n_estimators=250)"
$ws.Range("G13").Value = "significantly improves over polynomial regression "

# Row 12 H cell text changes (old "have to check with news..." note is replaced)
$ws.Range("H12").Value = "PolynomialFeatures between 2 and 6 (depending upon data length) works best"

# New row 14: xgboost approach
$ws.Range("C14").Value = "xgboost"
$ws.Range("D14").Value = "do"
$ws.Range("E14").Value = "NA"
$ws.Range("F14").Value = "do"
$ws.Range("G14").Value = "slightly improves over random forest"
$ws.Range("H14").Value = "the zoomed in fit is a lot better but train/test split doesn't work as good but can be parameterized better"

# Adjust column widths for columns I:J to match the new layout
$ws.Range("I1:J1").EntireColumn.ColumnWidth = 8.5

# Update selection to match the author's final cursor position
$ws.Range("L5").Select()
